# Insert two new data rows at position 110 (pushes old rows 110-199 down to 112-201)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("110:111").Insert()

# New row 110: Choclero / Primera, fecha 2022-03-31, volumen 8000
$ws.Cells.Item(110,1).Value2  = 7
$ws.Cells.Item(110,2).Value2  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(110,3).Value2  = "Ñuble"
$ws.Cells.Item(110,4).Value2  = 44651
$ws.Cells.Item(110,5).Value2  = 16
$ws.Cells.Item(110,6).Value2  = 100112024
$ws.Cells.Item(110,7).Value2  = "Choclo"
$ws.Cells.Item(110,8).Value2  = "Choclero"
$ws.Cells.Item(110,9).Value2  = "Primera"
$ws.Cells.Item(110,10).Value2 = 8000
$ws.Cells.Item(110,11).Value2 = 200
$ws.Cells.Item(110,12).Value2 = 200
$ws.Cells.Item(110,13).Value2 = 200
$ws.Cells.Item(110,14).Value2 = "`$/unidad"
$ws.Cells.Item(110,15).Value2 = "Región del Maule"
$ws.Cells.Item(110,16).Value2 = 200
$ws.Cells.Item(110,17).Value2 = 1
$ws.Cells.Item(110,18).Value2 = "Hortaliza"

# New row 111: Choclero / Segunda, fecha 2022-03-31, volumen 8000
$ws.Cells.Item(111,1).Value2  = 7
$ws.Cells.Item(111,2).Value2  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(111,3).Value2  = "Ñuble"
$ws.Cells.Item(111,4).Value2  = 44651
$ws.Cells.Item(111,5).Value2  = 16
$ws.Cells.Item(111,6).Value2  = 100112024
$ws.Cells.Item(111,7).Value2  = "Choclo"
$ws.Cells.Item(111,8).Value2  = "Choclero"
$ws.Cells.Item(111,9).Value2  = "Segunda"
$ws.Cells.Item(111,10).Value2 = 8000
$ws.Cells.Item(111,11).Value2 = 150
$ws.Cells.Item(111,12).Value2 = 150
$ws.Cells.Item(111,13).Value2 = 150
$ws.Cells.Item(111,14).Value2 = "`$/unidad"
$ws.Cells.Item(111,15).Value2 = "Región del Maule"
$ws.Cells.Item(111,16).Value2 = 150
$ws.Cells.Item(111,17).Value2 = 1
$ws.Cells.Item(111,18).Value2 = "Hortaliza"
